$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 keeps the text label "HK_G_acc_SD" (new shared-string entry in the diff,
# but same visible text as before)
$ws.Range("A1").Value = "HK_G_acc_SD"

# New numeric results for the HK_G model (rows 2-49)
$ws.Cells.Item(2, 1).Value = 45.25745257452575
$ws.Cells.Item(3, 1).Value = 43.089430894308947
$ws.Cells.Item(4, 1).Value = 44.715447154471541
$ws.Cells.Item(5, 1).Value = 44.986449864498645
$ws.Cells.Item(6, 1).Value = 45.25745257452575
$ws.Cells.Item(7, 1).Value = 45.25745257452575
$ws.Cells.Item(8, 1).Value = 45.528455284552841
$ws.Cells.Item(9, 1).Value = 49.322493224932252
$ws.Cells.Item(10, 1).Value = 46.341463414634148
$ws.Cells.Item(11, 1).Value = 45.799457994579946
$ws.Cells.Item(12, 1).Value = 43.089430894308947
$ws.Cells.Item(13, 1).Value = 46.883468834688344
$ws.Cells.Item(14, 1).Value = 47.154471544715449
$ws.Cells.Item(15, 1).Value = 46.883468834688344
$ws.Cells.Item(16, 1).Value = 46.612466124661246
$ws.Cells.Item(17, 1).Value = 46.341463414634148
$ws.Cells.Item(18, 1).Value = 46.070460704607044
$ws.Cells.Item(19, 1).Value = 48.780487804878049
$ws.Cells.Item(20, 1).Value = 45.25745257452575
$ws.Cells.Item(21, 1).Value = 45.25745257452575
$ws.Cells.Item(22, 1).Value = 46.341463414634148
$ws.Cells.Item(23, 1).Value = 41.192411924119241
$ws.Cells.Item(24, 1).Value = 38.211382113821138
$ws.Cells.Item(25, 1).Value = 37.94037940379404
$ws.Cells.Item(26, 1).Value = 47.696476964769644
$ws.Cells.Item(27, 1).Value = 45.799457994579946
$ws.Cells.Item(28, 1).Value = 50.948509485094853
$ws.Cells.Item(29, 1).Value = 45.799457994579946
$ws.Cells.Item(30, 1).Value = 46.883468834688344
$ws.Cells.Item(31, 1).Value = 47.154471544715449
$ws.Cells.Item(32, 1).Value = 41.192411924119241
$ws.Cells.Item(33, 1).Value = 40.650406504065039
$ws.Cells.Item(34, 1).Value = 42.005420054200542
$ws.Cells.Item(35, 1).Value = 39.295392953929536
$ws.Cells.Item(36, 1).Value = 37.669376693766935
$ws.Cells.Item(37, 1).Value = 44.986449864498645
$ws.Cells.Item(38, 1).Value = 37.398373983739837
$ws.Cells.Item(39, 1).Value = 37.94037940379404
$ws.Cells.Item(40, 1).Value = 38.482384823848236
$ws.Cells.Item(41, 1).Value = 47.154471544715449
$ws.Cells.Item(42, 1).Value = 47.154471544715449
$ws.Cells.Item(43, 1).Value = 47.154471544715449
$ws.Cells.Item(44, 1).Value = 45.25745257452575
$ws.Cells.Item(45, 1).Value = 45.25745257452575
$ws.Cells.Item(46, 1).Value = 44.715447154471541
$ws.Cells.Item(47, 1).Value = 43.360433604336045
$ws.Cells.Item(48, 1).Value = 49.322493224932252
$ws.Cells.Item(49, 1).Value = 45.799457994579946
